# Weekly update: a new price-report row is inserted at row 199 (pushing the
# existing rows 199-261 down to 200-262), matching the commit
# "Fruta / hortaliza, semanal" (fruit/vegetable, weekly refresh).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 199; everything below shifts down one row.
$ws.Rows.Item(199).Insert()

# Populate the newly inserted row 199 with this week's data point.
$ws.Cells.Item(199, 1).Value = 8
$ws.Cells.Item(199, 2).Value = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(199, 3).Value = 'Coquimbo'
$ws.Cells.Item(199, 4).Value = 44900
$ws.Cells.Item(199, 5).Value = 4
$ws.Cells.Item(199, 6).Value = 100112037
$ws.Cells.Item(199, 7).Value = 'Cebollín'
$ws.Cells.Item(199, 8).Value = 'Sin especificar'
$ws.Cells.Item(199, 9).Value = 'Primera'
$ws.Cells.Item(199, 10).Value = 600
$ws.Cells.Item(199, 11).Value = 5500
$ws.Cells.Item(199, 12).Value = 6000
$ws.Cells.Item(199, 13).Value = 5750
$ws.Cells.Item(199, 14).Value = '$/paquete 36 unidades'
$ws.Cells.Item(199, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(199, 16).Value = 160
$ws.Cells.Item(199, 17).Value = 36
$ws.Cells.Item(199, 18).Value = 'Hortaliza'
